# Issue #68: Implement smart default behavior for Advanced Alerts
# -----------------------------------------------------------------
# This script reproduces the author's editing session:
#  - tidy up the "2" placeholder sentinel values left in the demo
#    worksheets' first column down to the real default value "1"
#    on the SMS Tests / SMS Other Tests / Excel Number Tests sheets
#  - select the relevant range on "SMS Other Tests" (A2:A9) that the
#    author was reviewing
#  - leave "Simple Alerts" as the active tab when the workbook is
#    saved (instead of "Advanced Alerts")

$wb = $excel.ActiveWorkbook

# --- SMS Tests: column A default sentinel 2 -> 1 (rows 2-13) ---
$wsSms = $wb.Worksheets.Item("SMS Tests")
$wsSms.Range("A2:A13").Value = 1

# --- SMS Other Tests: column A default sentinel 2 -> 1 (rows 2-9) ---
$wsSmsOther = $wb.Worksheets.Item("SMS Other Tests")
$wsSmsOther.Range("A2:A9").Value = 1

# --- Excel Number Tests: column A default sentinel 2 -> 1 (rows 2-3) ---
$wsExcelNum = $wb.Worksheets.Item("Excel Number Tests")
$wsExcelNum.Range("A2:A3").Value = 1

# --- leave the reviewed range selected on "SMS Other Tests" ---
$wsSmsOther.Range("A2:A9").Select()

# --- make "Simple Alerts" the active/selected tab on save ---
$wsSimple = $wb.Worksheets.Item("Simple Alerts")
$wsSimple.Activate()
